$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 19, shifting existing rows 19-26 down to 20-27
$ws.Rows.Item(19).Insert()

# Populate the newly inserted row 19 with the new data record
$ws.Cells.Item(19, 1).Value = 11
$ws.Cells.Item(19, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(19, 3).Value = "Bíobío"
$ws.Cells.Item(19, 4).Value = 44673
$ws.Cells.Item(19, 5).Value = 8
$ws.Cells.Item(19, 6).Value = 100112026
$ws.Cells.Item(19, 7).Value = "Haba"
$ws.Cells.Item(19, 8).Value = "Sin especificar"
$ws.Cells.Item(19, 9).Value = "Primera"
$ws.Cells.Item(19, 10).Value = 80
$ws.Cells.Item(19, 11).Value = 18000
$ws.Cells.Item(19, 12).Value = 19000
$ws.Cells.Item(19, 13).Value = 18375
$ws.Cells.Item(19, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(19, 15).Value = "Región Metropolitana"
$ws.Cells.Item(19, 16).Value = 735
$ws.Cells.Item(19, 17).Value = 25
$ws.Cells.Item(19, 18).Value = "Hortaliza"
